$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("English", 21.59770217433532),
  @("Chinese", 18.54551662976611),
  @("Spanish", 6.59276099873834),
  @("German", 4.30372386983408),
  @("Arabic", 4.247115533310451),
  @("Japanese", 4.03996908203427),
  @("Russian", 3.220604241080689),
  @("Malay-Indonesian", 3.151907800318338),
  @("Portuguese", 2.861751308099515),
  @("French", 2.619545008147039),
  @("Italian", 1.969904988929833),
  @("Turkish", 1.818988959772109),
  @("Korean", 1.701512093732438),
  @("Dutch", 1.228598994538466),
  @("Persian", 0.9965705998917499),
  @("Thai", 0.9854954241269278),
  @("Polish", 0.9698233947787287),
  @("Urdu", 0.8892694291320762),
  @("Vietnamese", 0.7540202725957748),
  @("Bengali", 0.7281202440245618)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove now-unused rows 22 and 23 (previously Uzbek and Vietnamese), shifting cells up
$ws.Range("A22:B23").Delete() | Out-Null
